$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper pattern used throughout: find a literal phrase, Delete() the
# range and InsertAfter() the replacement text rather than just
# overwriting Range.Text. Plain Range.Text / Find-replace keeps the
# matched run's direct character formatting (including w:shd highlight
# shading), which we do NOT want on the new templated text. Deleting
# the range first and inserting fresh text drops that inherited
# formatting, and we then explicitly (re)apply Font.Name / Font.Bold
# to match the target look.
# ---------------------------------------------------------------------

# 1) Date line: "20 dias del mes de mayo del ano 2025. -" -> "{{ fecha }} -"
$rng = $d.Content
$found = $rng.Find.Execute("20 días del mes de mayo del año 2025. -")
if ($found) {
    $rng.Delete()
    $rng.InsertAfter("{{ fecha }} -")
    $rng.Font.Name = "Georgia"
    $rng.Font.Bold = 0
}

# 2) Name + C.I. -> "{{ nombre_completo }} con C.I Nº {{ ci }}" (bold)
$rng = $d.Content
$found = $rng.Find.Execute("GABRIEL BAEZ SANCHEZ con C.I Nº 5.987.220")
if ($found) {
    $rng.Delete()
    $rng.InsertAfter("{{ nombre_completo }} con C.I Nº {{ ci }}")
    $rng.Font.Name = "Georgia"
    $rng.Font.Bold = 1
}

# 3) Street address -> "{{ dirección_calle }} de la Ciudad de {{ ciudad }}"
#    (the Find phrase also swallows the leftover " " connector run that
#    sits between the bold C.I. block and "con domicilio..." so that
#    run's stray explicit sz/szCs does not leak into the new text)
$rng = $d.Content
$found = $rng.Find.Execute(" con domicilio en las calles Calle Avda. San Francisco de la Ciudad de Limpio,")
if ($found) {
    $rng.Delete()
    $rng.InsertAfter(" con domicilio en las calles {{ dirección_calle }} de la Ciudad de {{ ciudad }},")
    $rng.Font.Name = "Georgia"
    $rng.Font.Bold = 0
}

# 4) Company + RUC -> "{{ empresa_que_trabajo }}, con RUC Nº {{ ruc_empresa }}" (bold) + plain " -"
$rng = $d.Content
$found = $rng.Find.Execute("CRISTHIAN SANTO TOMAS AQUINO AQUINO, con RUC Nº 3446196-5 –")
if ($found) {
    $rng.Delete()
    $rng.InsertAfter("{{ empresa_que_trabajo }}, con RUC Nº {{ ruc_empresa }}")
    $rng.Font.Name = "Georgia"
    $rng.Font.Bold = 1
    $rng.Collapse(0)
    $rng.InsertAfter(" –")
    $rng.Font.Name = "Georgia"
    $rng.Font.Bold = 0
}

# 5) Drop the stale lastRenderedPageBreak cached on "C.I. N°:" by
#    rewriting that run from scratch.
$rng = $d.Content
$found = $rng.Find.Execute("C.I. N°:")
if ($found) {
    $rng.Delete()
    $rng.InsertAfter("C.I. N°:")
    $rng.Font.Name = "Georgia"
    $rng.Font.Bold = 1
    $rng.Font.Size = 12
}
